$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citywide Totals")
Write-Host $ws.Name
$ws.Range("D2").Value = 95
